# Applies the "Deploying to gh-pages" FHIR StructureDefinition metadata
# refresh (admit-count, 5.0.0 -> 6.0.0) to the workbook.
#
# Sheet 1 ("Metadata", a Property/Value table):
#   - Version bumps from 5.0.0 to 6.0.0
#   - Date bumps to the new build timestamp
#   - Publisher gains a display value ("Alvearie Team")
#   - The duplicated "Contact" row is replaced by a single "Jurisdiction"
#     row ("United States of America"), which nets out to one fewer row
#     overall, so the second duplicate row is removed entirely.
#
# Sheet 2 ("Elements", the FHIR element table):
#   - The root Extension row's Short/Definition columns (K2/L2) are
#     updated to describe this specific extension instead of the generic
#     placeholder text.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item(1)

# Remove the second, duplicate "Contact" row (old row 11) - this shifts
# every following row up by one and shrinks the sheet from 21 to 20 rows.
$wsMeta.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value = "6.0.0"

# Date: refreshed build timestamp
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$wsMeta.Range("B9").Value = "Alvearie Team"

# What used to be the first "Contact" row becomes "Jurisdiction"
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

$wsElements = $wb.Worksheets.Item(2)

# Root Extension element: Short / Definition text
$wsElements.Range("K2").Value = "Admit Count"
$wsElements.Range("L2").Value = "Used to determine which facility claims should be counted as inpatient admissions"
